$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C3").Value = [double]"1"
$ws.Range("E3").Value = [double]"0"
$ws.Range("C4").Value = [double]"0.92777620480442979"
$ws.Range("E4").Value = [double]"0.072223795195570262"
$ws.Range("C5").Value = [double]"0.73156848902157257"
$ws.Range("E5").Value = [double]"0.00075404922037147759"
$ws.Range("C6").Value = [double]"0.44860486960559709"
$ws.Range("E6").Value = [double]"0.068412932907872953"
$ws.Range("C7").Value = [double]"0.15293852533765273"
$ws.Range("E7").Value = [double]"0.028139083648352728"
$ws.Range("C8").Value = [double]"0.075637625912831158"
$ws.Range("E8").Value = [double]"0.0283569391067643"
$ws.Range("C9").Value = [double]"0.031775865014838665"
$ws.Range("E9").Value = [double]"0.0090074153674495686"
$ws.Range("C10").Value = [double]"0.013374580219945247"
$ws.Range("E10").Value = [double]"0.0049454214429347004"
$ws.Range("C11").Value = [double]"0.0040008204008844727"
$ws.Range("E11").Value = [double]"5.132577869039652e-05"
$ws.Range("C13").Value = [double]"1"
$ws.Range("E13").Value = [double]"0"
$ws.Range("C14").Value = [double]"1.0199387227673973"
$ws.Range("E14").Value = [double]"0.019938722767397277"
$ws.Range("C15").Value = [double]"0.91365923851480035"
$ws.Range("E15").Value = [double]"0.015111640714474073"
$ws.Range("C16").Value = [double]"0.65293828229087914"
$ws.Range("E16").Value = [double]"0.059831743802541915"
$ws.Range("C17").Value = [double]"0.36654835753681869"
$ws.Range("E17").Value = [double]"0.035475837575819756"
$ws.Range("C18").Value = [double]"0.14293494630484724"
$ws.Range("E18").Value = [double]"0.026247068951715625"
$ws.Range("C19").Value = [double]"0.071827198854823721"
$ws.Range("E19").Value = [double]"0.00037473645294306546"
$ws.Range("C20").Value = [double]"0.033856757712352505"
$ws.Range("E20").Value = [double]"0.0028455821051713518"
$ws.Range("C21").Value = [double]"0.028737813790577531"
$ws.Range("E21").Value = [double]"0.00092185323473713266"
$ws.Range("C23").Value = [double]"1"
$ws.Range("E23").Value = [double]"0"
$ws.Range("C24").Value = [double]"0.99721976734366358"
$ws.Range("E24").Value = [double]"0.0027802326563364783"
$ws.Range("C25").Value = [double]"0.6874565028611046"
$ws.Range("E25").Value = [double]"0.065046322737951245"
$ws.Range("C26").Value = [double]"0.26608835331055697"
$ws.Range("E26").Value = [double]"0.021154206231468062"
$ws.Range("C27").Value = [double]"0.089213220362641635"
$ws.Range("E27").Value = [double]"0.027376813995345396"
$ws.Range("C28").Value = [double]"0.025202537218857138"
$ws.Range("E28").Value = [double]"0.011295308279021586"
$ws.Range("C29").Value = [double]"0.011182587480671271"
$ws.Range("E29").Value = [double]"0.0032609819874478874"
$ws.Range("C30").Value = [double]"0.016272482841773687"
$ws.Range("E30").Value = [double]"0.008154636712533643"
$ws.Range("C31").Value = [double]"0.0043045535164640131"
$ws.Range("E31").Value = [double]"0.00010259988803462726"
$ws.Range("C33").Value = [double]"1"
$ws.Range("E33").Value = [double]"0"
$ws.Range("C34").Value = [double]"0.95869972783601165"
$ws.Range("E34").Value = [double]"0.041300272163988405"
$ws.Range("C35").Value = [double]"0.9078588854278602"
$ws.Range("E35").Value = [double]"0.048662646946527943"
$ws.Range("C36").Value = [double]"0.56957005607573574"
$ws.Range("E36").Value = [double]"0.0063594554874340506"
$ws.Range("C37").Value = [double]"0.23720313823263667"
$ws.Range("E37").Value = [double]"0.012238919802421241"
$ws.Range("C38").Value = [double]"0.073704579902514791"
$ws.Range("E38").Value = [double]"0.01480073235646898"
$ws.Range("C39").Value = [double]"0.06629753361429544"
$ws.Range("E39").Value = [double]"0.0048097104061536541"
$ws.Range("C40").Value = [double]"0.050795534750370137"
$ws.Range("E40").Value = [double]"0.0096681249227956439"
$ws.Range("C41").Value = [double]"0.02864383112886227"
$ws.Range("E41").Value = [double]"0.00066393072138239476"
Write-Output "done"
